$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Report header label edits (surgical in-run text replace) ---
# si 6: "Volume 32   Number  12" -> "...13" (last run "12" -> "13")
$ws.Range("A8").Characters(21, 2).Text = "13"
# si 9: week-of dates "3/17/2025".."3/23/2025" -> "3/24/2025".."3/30/2025"
$ws.Range("C9").Characters(27, 9).Text = "3/24/2025"
$ws.Range("C9").Characters(47, 9).Text = "3/30/2025"

# --- Updated weekly crime statistics (rows 14-33) ---
# C31 (Hate Crimes, Week-to-Date '25) goes from a blank/zero text cell to a
# real numeric count, so it picks up the same numeric style as its column.
$ws.Range("C31").NumberFormat = "#,##0"
$ws.Range("D14").Value2 = 4
$ws.Range("G14").Value2 = 6
$ws.Range("H14").Value2 = -83.333333333333
$ws.Range("J14").Value2 = 15
$ws.Range("K14").Value2 = -66.666666666666
$ws.Range("L14").Value2 = -66.666666666666
$ws.Range("M14").Value2 = -58.333333333333
$ws.Range("N14").Value2 = -93.902439024390
$ws.Range("C15").Value2 = 4
$ws.Range("D15").Value2 = 2
$ws.Range("E15").Value2 = 100
$ws.Range("I15").Value2 = 53
$ws.Range("J15").Value2 = 29
$ws.Range("K15").Value2 = 82.758620689655
$ws.Range("L15").Value2 = 51.428571428571
$ws.Range("M15").Value2 = 15.217391304347
$ws.Range("N15").Value2 = -50.467289719626
$ws.Range("C16").Value2 = 31
$ws.Range("D16").Value2 = 48
$ws.Range("E16").Value2 = -35.416666666666
$ws.Range("F16").Value2 = 122
$ws.Range("G16").Value2 = 164
$ws.Range("H16").Value2 = -25.609756097561
$ws.Range("I16").Value2 = 400
$ws.Range("J16").Value2 = 498
$ws.Range("K16").Value2 = -19.678714859437
$ws.Range("L16").Value2 = -7.407407407407
$ws.Range("M16").Value2 = -26.199261992619
$ws.Range("N16").Value2 = -81.924988703117
$ws.Range("C17").Value2 = 57
$ws.Range("D17").Value2 = 53
$ws.Range("E17").Value2 = 7.547169811320
$ws.Range("F17").Value2 = 204
$ws.Range("G17").Value2 = 256
$ws.Range("H17").Value2 = -20.3125
$ws.Range("I17").Value2 = 623
$ws.Range("J17").Value2 = 712
$ws.Range("K17").Value2 = -12.5
$ws.Range("L17").Value2 = -5.319148936170
$ws.Range("M17").Value2 = 50.120481927710
$ws.Range("N17").Value2 = -48.169717138103
$ws.Range("F18").Value2 = 102
$ws.Range("G18").Value2 = 114
$ws.Range("H18").Value2 = -10.526315789473
$ws.Range("I18").Value2 = 322
$ws.Range("J18").Value2 = 358
$ws.Range("K18").Value2 = -10.055865921787
$ws.Range("L18").Value2 = -23.333333333333
$ws.Range("M18").Value2 = -7.204610951008
$ws.Range("N18").Value2 = -87.770603873908
$ws.Range("C19").Value2 = 97
$ws.Range("D19").Value2 = 105
$ws.Range("E19").Value2 = -7.619047619047
$ws.Range("F19").Value2 = 466
$ws.Range("G19").Value2 = 487
$ws.Range("H19").Value2 = -4.312114989733
$ws.Range("I19").Value2 = 1304
$ws.Range("J19").Value2 = 1502
$ws.Range("K19").Value2 = -13.182423435419
$ws.Range("L19").Value2 = -6.790564689063
$ws.Range("M19").Value2 = 25.505293551491
$ws.Range("N19").Value2 = -49.260700389105
$ws.Range("C20").Value2 = 17
$ws.Range("D20").Value2 = 20
$ws.Range("E20").Value2 = -15
$ws.Range("F20").Value2 = 68
$ws.Range("G20").Value2 = 75
$ws.Range("H20").Value2 = -9.333333333333
$ws.Range("I20").Value2 = 170
$ws.Range("J20").Value2 = 222
$ws.Range("K20").Value2 = -23.423423423423
$ws.Range("L20").Value2 = -39.716312056737
$ws.Range("M20").Value2 = 70
$ws.Range("N20").Value2 = -92.747440273037
$ws.Range("C21").Value2 = 232
$ws.Range("D21").Value2 = 259
$ws.Range("E21").Value2 = -10.424710424710
$ws.Range("F21").Value2 = 978
$ws.Range("G21").Value2 = 1107
$ws.Range("H21").Value2 = -11.653116531165
$ws.Range("I21").Value2 = 2877
$ws.Range("J21").Value2 = 3336
$ws.Range("K21").Value2 = -13.758992805755
$ws.Range("L21").Value2 = -11.231101511879
$ws.Range("M21").Value2 = 15.033986405437
$ws.Range("N21").Value2 = -74.199623352165
$ws.Range("C22").Value2 = 5
$ws.Range("D22").Value2 = 6
$ws.Range("E22").Value2 = -16.666666666666
$ws.Range("F22").Value2 = 18
$ws.Range("G22").Value2 = 20
$ws.Range("H22").Value2 = -10
$ws.Range("I22").Value2 = 56
$ws.Range("J22").Value2 = 74
$ws.Range("K22").Value2 = -24.324324324324
$ws.Range("L22").Value2 = -20
$ws.Range("M22").Value2 = 0
$ws.Range("C23").Value2 = 24
$ws.Range("D23").Value2 = 27
$ws.Range("E23").Value2 = -11.111111111111
$ws.Range("F23").Value2 = 95
$ws.Range("G23").Value2 = 97
$ws.Range("H23").Value2 = -2.061855670103
$ws.Range("I23").Value2 = 301
$ws.Range("J23").Value2 = 316
$ws.Range("K23").Value2 = -4.746835443037
$ws.Range("L23").Value2 = 3.793103448275
$ws.Range("M23").Value2 = 54.358974358974
$ws.Range("C24").Value2 = 323
$ws.Range("D24").Value2 = 259
$ws.Range("E24").Value2 = 24.710424710424
$ws.Range("F24").Value2 = 1124
$ws.Range("G24").Value2 = 965
$ws.Range("H24").Value2 = 16.476683937823
$ws.Range("I24").Value2 = 3506
$ws.Range("J24").Value2 = 2963
$ws.Range("K24").Value2 = 18.326020924738
$ws.Range("L24").Value2 = 6.500607533414
$ws.Range("M24").Value2 = 77.969543147208
$ws.Range("C25").Value2 = 170
$ws.Range("D25").Value2 = 157
$ws.Range("E25").Value2 = 8.280254777070
$ws.Range("F25").Value2 = 651
$ws.Range("G25").Value2 = 505
$ws.Range("H25").Value2 = 28.910891089108
$ws.Range("I25").Value2 = 2002
$ws.Range("J25").Value2 = 1576
$ws.Range("K25").Value2 = 27.030456852791
$ws.Range("L25").Value2 = 7.345844504021
$ws.Range("C26").Value2 = 100
$ws.Range("D26").Value2 = 95
$ws.Range("E26").Value2 = 5.263157894736
$ws.Range("F26").Value2 = 375
$ws.Range("G26").Value2 = 392
$ws.Range("H26").Value2 = -4.336734693877
$ws.Range("I26").Value2 = 1063
$ws.Range("J26").Value2 = 1112
$ws.Range("K26").Value2 = -4.406474820143
$ws.Range("L26").Value2 = 3.505355404089
$ws.Range("M26").Value2 = -9.762308998302
$ws.Range("C27").Value2 = 5
$ws.Range("D27").Value2 = 4
$ws.Range("E27").Value2 = 25
$ws.Range("I27").Value2 = 63
$ws.Range("J27").Value2 = 48
$ws.Range("K27").Value2 = 31.25
$ws.Range("L27").Value2 = 1.612903225806
$ws.Range("C28").Value2 = 8
$ws.Range("D28").Value2 = 8
$ws.Range("E28").Value2 = 0
$ws.Range("F28").Value2 = 52
$ws.Range("G28").Value2 = 48
$ws.Range("H28").Value2 = 8.333333333333
$ws.Range("I28").Value2 = 132
$ws.Range("J28").Value2 = 135
$ws.Range("K28").Value2 = -2.222222222222
$ws.Range("L28").Value2 = -4.347826086956
$ws.Range("C29").Value2 = 3
$ws.Range("D29").Value2 = 2
$ws.Range("E29").Value2 = 50
$ws.Range("F29").Value2 = 7
$ws.Range("H29").Value2 = 0
$ws.Range("I29").Value2 = 17
$ws.Range("J29").Value2 = 23
$ws.Range("K29").Value2 = -26.086956521739
$ws.Range("L29").Value2 = -54.054054054054
$ws.Range("M29").Value2 = -51.428571428571
$ws.Range("N29").Value2 = -91.237113402061
$ws.Range("G30").Value2 = 5
$ws.Range("H30").Value2 = 20
$ws.Range("I30").Value2 = 15
$ws.Range("J30").Value2 = 18
$ws.Range("K30").Value2 = -16.666666666666
$ws.Range("L30").Value2 = -55.882352941176
$ws.Range("M30").Value2 = -54.545454545454
$ws.Range("N30").Value2 = -91.620111731843
$ws.Range("C31").Value2 = 3
$ws.Range("D31").Value2 = 3
$ws.Range("E31").Value2 = 0
$ws.Range("F31").Value2 = 7
$ws.Range("H31").Value2 = -41.666666666666
$ws.Range("I31").Value2 = 16
$ws.Range("J31").Value2 = 26
$ws.Range("K31").Value2 = -38.461538461538
$ws.Range("L31").Value2 = -11.111111111111
$ws.Range("J33").Value2 = 3
$ws.Range("K33").Value2 = 33.333333333333
